$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Establish shared-string creation order to match target workbook:
# 16 "Entry ", 17 "Options File failed to open", 18 "Options file has no window name",
# 19 "Options file has no window size", 20 "Options file window size formatted incorrectly",
# 21 "Options file has no on load nut", 22 "main"

# 1. Create "Entry " first (index 16) via D13
$ws.Range("D13").Value = "Entry "

# 2-6. Create the new "Meaning" strings for rows 13-17 (indices 17-21)
$ws.Range("C13").Value = "Options File failed to open"
$ws.Range("C14").Value = "Options file has no window name"
$ws.Range("C15").Value = "Options file has no window size"
$ws.Range("C16").Value = "Options file window size formatted incorrectly"
$ws.Range("C17").Value = "Options file has no on load nut"

# 7. Fill remaining "Entry " cells (reuse existing shared string)
$ws.Range("D14").Value = "Entry "
$ws.Range("D15").Value = "Entry "
$ws.Range("D16").Value = "Entry "
$ws.Range("D17").Value = "Entry "

# 8. Create "main" last (index 22) and fill E13:E17
$ws.Range("E13").Value = "main"
$ws.Range("E14").Value = "main"
$ws.Range("E15").Value = "main"
$ws.Range("E16").Value = "main"
$ws.Range("E17").Value = "main"

# Update selection to match target active cell
$ws.Range("G15").Select()
